$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.213.99'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.911.51'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''0.7385'
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").Value = '''243.71'
$ws.Range("E6").Value = '  -2.09%  '
$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.3135'
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").Value = '''26.99'
$ws.Range("E9").Value = '  -4.40%  '
$ws.Range("D10").Value = '''0.06955'
$ws.Range("E10").Value = '  -2.23%  '
$ws.Range("D11").Value = '''0.7774'
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").Value = '''0.07971'
$ws.Range("D13").Value = '1.965.79'
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").Value = '''5.270'
$ws.Range("E14").Value = '  -2.15%  '
$ws.Range("D15").Value = '''91.37'
$ws.Range("E15").Value = '  -3.34%  '
$ws.Range("D16").Value = '30.252.41'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '''14.26'
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("D18").Value = '''244.64'
$ws.Range("E18").Value = '  -3.39%  '
$ws.Range("D19").Value = '''5.803'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = '''0.000007812'
$ws.Range("E20").Value = '  -2.82%  '
$ws.Range("D21").Value = '2.182.98'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '''6.626'
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("D25").Value = '''9.394'
$ws.Range("E25").Value = '  -1.84%  '
$ws.Range("D26").Value = '''165.13'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = '''19.01'
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("D28").Value = '''0.1268'
$ws.Range("E28").Value = '  -4.39%  '
$ws.Range("D29").Value = '''2.129'
$ws.Range("E29").Value = '  -8.93%  '
$ws.Range("D30").Value = '''1.353'
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").Value = '''1.547'
$ws.Range("E31").Value = '  +1.00%  '
$ws.Range("D32").Value = '''4.313'
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").Value = '''4.083'
$ws.Range("E33").Value = '  -1.60%  '
$ws.Range("D34").Value = '''0.05201'
$ws.Range("E34").Value = '  +1.34%  '
$ws.Range("D35").Value = '''1.294'
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("D36").Value = '''0.7516'
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").Value = '''2.761'
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("D38").Value = '''0.01939'
$ws.Range("E38").Value = '  -1.45%  '
$ws.Range("D39").Value = '''2.792'
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("D40").Value = '''6.407'
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").Value = '''75.96'
$ws.Range("E41").Value = '  -2.87%  '
$ws.Range("D42").Value = '''0.4483'
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("D43").Value = '''1.947'
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").Value = '''1.000'
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").Value = '''0.8327'
$ws.Range("E45").Value = '  -1.27%  '
$ws.Range("D46").Value = '''7.660'
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''9.920'
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''101.36'
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").Value = '2.061.59'
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").Value = '''37.03'
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("D51").Value = '''0.1215'
$ws.Range("E51").Value = '  +1.41%  '
